$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Gnai2", "P2ry12", "M1", 3, 1, 195.435389, 586.306167, 0.3095741734129938, 0.3095741734129938, 3, 1, 13.757307, 41.27192100000001, 0.3428896105196262, 0.3428896105196262, 2688.664645137423, 24197.98180623681, 0.1061497677485167, 0.1061497677485166),
  @("ECs", "Gnai2", "P2ry12", "M2", 3, 1, 195.435389, 586.306167, 0.3095741734129938, 0.3095741734129938, 3, 1, 26.34504533333333, 79.03513599999999, 0.6566286798330924, 0.6566286798330924, 5148.754182942634, 46338.7876464837, 0.2032752807985949, 0.2032752807985949),
  @("ECs", "Gnai2", "P2ry12", "sCs", 3, 1, 195.435389, 586.306167, 0.3095741734129938, 0.3095741734129938, 1, 0.3333333333333333, 0.019327, 0.057981, 0.0004817096472814639, 0.0004817096472814639, 3.777179763203, 33.994617868827, 0.000149124865882224, 0.000149124865882224),
  @("FAPs", "Gnai2", "P2ry12", "M1", 3, 1, 79.82725266666667, 239.481758, 0.1264482133280045, 0.1264482133280045, 3, 1, 13.757307, 41.27192100000001, 0.3428896105196262, 0.3428896105196262, 1098.208021901902, 9883.872197117118, 0.04335777861894206, 0.04335777861894206),
  @("FAPs", "Gnai2", "P2ry12", "M2", 3, 1, 79.82725266666667, 239.481758, 0.1264482133280045, 0.1264482133280045, 3, 1, 26.34504533333333, 79.03513599999999, 0.6566286798330924, 0.6566286798330924, 2103.052590338787, 18927.47331304909, 0.08302952338482082, 0.08302952338482081),
  @("FAPs", "Gnai2", "P2ry12", "sCs", 3, 1, 79.82725266666667, 239.481758, 0.1264482133280045, 0.1264482133280045, 1, 0.3333333333333333, 0.019327, 0.057981, 0.0004817096472814639, 0.0004817096472814639, 1.542821312288667, 13.885391810598, 0.00006091132424160435, 0.00006091132424160434),
  @("M1", "Gnai2", "P2ry12", "M1", 3, 1, 162.5116576666667, 487.534973, 0.2574222219914007, 0.2574222219914007, 3, 1, 13.757307, 41.27192100000001, 0.3428896105196262, 0.3428896105196262, 2235.722765599237, 20121.50489039313, 0.08826740543772814, 0.08826740543772814),
  @("M1", "Gnai2", "P2ry12", "M2", 3, 1, 162.5116576666667, 487.534973, 0.2574222219914007, 0.2574222219914007, 3, 1, 26.34504533333333, 79.03513599999999, 0.6566286798330924, 0.6566286798330924, 4281.37698842348, 38532.39289581132, 0.1690308137859147, 0.1690308137859147),
  @("M1", "Gnai2", "P2ry12", "sCs", 3, 1, 162.5116576666667, 487.534973, 0.2574222219914007, 0.2574222219914007, 1, 0.3333333333333333, 0.019327, 0.057981, 0.0004817096472814639, 0.0004817096472814639, 3.140862807723666, 28.267765269513, 0.0001240027677578883, 0.0001240027677578883),
  @("M2", "Gnai2", "P2ry12", "M1", 3, 1, 159.7910413333334, 479.3731240000001, 0.2531127028358626, 0.2531127028358626, 3, 1, 13.757307, 41.27192100000001, 0.3428896105196262, 0.3428896105196262, 2198.294411472357, 19784.64970325121, 0.08678971609295882, 0.0867897160929588),
  @("M2", "Gnai2", "P2ry12", "M2", 3, 1, 159.7910413333334, 479.3731240000001, 0.2531127028358626, 0.2531127028358626, 3, 1, 26.34504533333333, 79.03513599999999, 0.6566286798330924, 0.6566286798330924, 4209.702227787207, 37887.32005008487, 0.1662010599120983, 0.1662010599120982),
  @("M2", "Gnai2", "P2ry12", "sCs", 3, 1, 159.7910413333334, 479.3731240000001, 0.2531127028358626, 0.2531127028358626, 1, 0.3333333333333333, 0.019327, 0.057981, 0.0004817096472814639, 0.0004817096472814639, 3.088281455849334, 27.794533102644, 0.0001219268308055214, 0.0001219268308055214),
  @("sCs", "Gnai2", "P2ry12", "M1", 3, 1, 33.73857866666666, 101.215736, 0.05344268843173843, 0.05344268843173842, 3, 1, 13.757307, 41.27192100000001, 0.3428896105196262, 0.3428896105196262, 464.151984460984, 4177.367860148856, 0.01832494262148052, 0.01832494262148052),
  @("sCs", "Gnai2", "P2ry12", "M2", 3, 1, 33.73857866666666, 101.215736, 0.05344268843173843, 0.05344268843173842, 3, 1, 26.34504533333333, 79.03513599999999, 0.6566286798330924, 0.6566286798330924, 888.844384455566, 7999.599460100095, 0.03509200195166368, 0.03509200195166368),
  @("sCs", "Gnai2", "P2ry12", "sCs", 3, 1, 33.73857866666666, 101.215736, 0.05344268843173843, 0.05344268843173842, 1, 0.3333333333333333, 0.019327, 0.057981, 0.0004817096472814639, 0.0004817096472814639, 0.6520655098906666, 5.868589589016, 0.00002574385859422589, 0.00002574385859422589),
)
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $rowData = $data[$i]
  for ($j = 0; $j -lt $rowData.Count; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $rowData[$j]
  }
}
